# Simulated Wild Card round and logged it

$wb = $excel.ActiveWorkbook

# --- Add the new WR entry (T.Benjamin) to the "WR" sheet ---
$wr = $wb.Worksheets.Item("WR")

$wr.Cells.Item(9, 1).Value = "T.Benjamin"
$wr.Cells.Item(9, 2).Value = 0
$wr.Cells.Item(9, 3).Value = 0
$wr.Cells.Item(9, 4).Value = 0
$wr.Cells.Item(9, 5).Value = 0
$wr.Cells.Item(9, 6).Value = 0
$wr.Cells.Item(9, 7).Value = 0
$wr.Cells.Item(9, 8).Value = 0
$wr.Cells.Item(9, 9).Value = 0
$wr.Cells.Item(9, 10).Value = 0

# Update selection on the WR sheet to the newly added row/range
$wr.Range("B8").Select()
$wr.Range("B8:J9").Select()

# WR is now the active/selected tab
$wr.Select()

# --- RB sheet is no longer the selected tab ---
$rb = $wb.Worksheets.Item("RB")
$rb.Range("K7").Select()

# Re-activate WR as the active sheet/tab last so it ends up selected
$wr.Activate()
